$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (1_control_dl) - values updated, test name unchanged
$ws.Range("E2").Value = 24.28000000000036
$ws.Range("G2").Value = [double]"1.754152378907747e-14"
$ws.Range("H2").Value = [double]"7.317117550897406e-14"
$ws.Range("K2").Value = 535.2267947162426
$ws.Range("L2").Value = "[377.68545930394953, 692.7681301285356]"
$ws.Range("M2").Value = [double]"6.684144349122789e-10"
$ws.Range("N2").Value = [double]"8.912192465497052e-10"
$ws.Range("O2").Value = 2.144710900648889
$ws.Range("P2").Value = "[1.8679740102425804, 2.421447791055197]"
$ws.Range("S2").Value = 781.6150962990639
$ws.Range("T2").Value = "[699.5290729246699, 863.7011196734579]"
$ws.Range("W2").Value = 15.99223223223246
$ws.Range("X2").Value = 14.92284284284306
$ws.Range("Y2").Value = 17.06162162162187

# Row 3 (2_control_dd -> 2_induction_dd)
$ws.Range("C3").Value = "2_induction_dd"
$ws.Range("E3").Value = 23.86000000000029
$ws.Range("H3").Value = [double]"4.989766402809693e-16"
$ws.Range("K3").Value = 249.4754662368658
$ws.Range("L3").Value = "[200.29626087150325, 298.6546716022283]"
$ws.Range("O3").Value = 2.484342538874811
$ws.Range("P3").Value = "[2.283079345852041, 2.6856057318975814]"
$ws.Range("S3").Value = 731.4380143936773
$ws.Range("T3").Value = "[705.7761632240703, 757.0998655632844]"
$ws.Range("W3").Value = 14.42586586586604
$ws.Range("X3").Value = 13.66158158158175
$ws.Range("Y3").Value = 15.19015015015034

# Row 4 (3_induction_dd -> 3_hypo_dd)
$ws.Range("C4").Value = "3_hypo_dd"
$ws.Range("E4").Value = 23.74000000000027
$ws.Range("G4").Value = [double]"1.788461992946466e-07"
$ws.Range("H4").Value = [double]"3.502495947018784e-07"
$ws.Range("K4").Value = 125.2990594716739
$ws.Range("L4").Value = "[71.77356901497711, 178.8245499283707]"
$ws.Range("M4").Value = [double]"5.426034240807454e-06"
$ws.Range("N4").Value = [double]"5.426034240807454e-06"
$ws.Range("O4").Value = 1.239026532046426
$ws.Range("P4").Value = "[0.735868549489501, 1.74218451460335]"
$ws.Range("Q4").Value = [double]"1.767423148946179e-06"
$ws.Range("R4").Value = [double]"1.767423148946179e-06"
$ws.Range("S4").Value = 743.4172392807579
$ws.Range("T4").Value = "[712.6704207473376, 774.1640578141782]"
$ws.Range("W4").Value = 19.05853853853876
$ws.Range("X4").Value = 17.15743743743763
$ws.Range("Y4").Value = 20.95963963963988

# Row 5 (4_hypo_dd -> 4_hypo_dl)
$ws.Range("C5").Value = "4_hypo_dl"
$ws.Range("E5").Value = 24.14000000000033
$ws.Range("H5").Value = [double]"4.989766402809693e-16"
$ws.Range("K5").Value = 450.7269942963157
$ws.Range("L5").Value = "[392.0594222221247, 509.3945663705067]"
$ws.Range("O5").Value = -2.603842559732081
$ws.Range("P5").Value = "[-2.7296320553713116, -2.47805306409285]"
$ws.Range("S5").Value = 475.8716733992901
$ws.Range("T5").Value = "[445.1047416910316, 506.6386051075486]"
$ws.Range("W5").Value = 10.0039639639641
$ws.Range("X5").Value = 9.520680680680815
$ws.Range("Y5").Value = 10.48724724724739

# Delete row 6 entirely (was 5_hypo_dl)
$ws.Rows.Item(6).Delete()
